# Updated symbol list on Sat Jan 28 21:50:20 UTC 2023 with GitHub Actions
# Refresh coin prices / volumes, and re-sort a block of exchange-token rows
# (rows 7-17) to reflect the new ranking order.
#
# Numeric-looking values (price/volume columns D & E) are written with a
# leading apostrophe so Excel keeps them as literal text (preserving
# trailing zeros / the "%" suffix) instead of coercing them into numbers,
# then the style is reset back to "Normal" so no stray text-format style
# id gets attached to the cell (matching the original workbook, which has
# no explicit style on these data cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.62%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'38.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'6.96%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.109"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.18%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08077"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.42%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.932"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-4.51%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.046"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.35%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9272"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.00%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1452"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-2.09%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1912"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.64%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09058"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.46%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03506"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.45%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09774"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-1.14%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001391"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.82%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005920"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.35%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.773"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.68%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.204"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.88%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-0.81%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.20%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1328"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.91%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.697"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.36%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2418"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.21%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04378"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.14%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001231"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.24%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004275"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.04%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.05%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02035"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-1.12%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05050"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.73%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007528"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.77%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009711"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.31%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-2.17%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-0.89%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009916"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.25%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006198"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.62%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.13%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002874"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.001804"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.13%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.13%"
$ws.Range("E51").Style = "Normal"
